# Auto-generated edit script
# Commit: Add data for 2024-02-16
# Updates 2024 partial-year violent crime totals (column K) across
# Citywide Totals, By Neighborhood, and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 845
$ws.Range("K3").Value = 799
$ws.Range("I4").Value = 1783
$ws.Range("J4").Value = 1794
$ws.Range("K4").Value = 180
$ws.Range("K5").Value = 47
$ws.Range("K6").Value = 1119
$ws.Range("I7").Value = 26237
$ws.Range("J7").Value = 29252
$ws.Range("K7").Value = 2990

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 88
$ws.Range("K8").Value = 183
$ws.Range("K11").Value = 64
$ws.Range("K19").Value = 78
$ws.Range("K22").Value = 10
$ws.Range("K23").Value = 28
$ws.Range("K25").Value = 16
$ws.Range("K27").Value = 36
$ws.Range("K29").Value = 149
$ws.Range("K31").Value = 32
$ws.Range("K33").Value = 129
$ws.Range("K39").Value = 5
$ws.Range("K41").Value = 28
$ws.Range("J42").Value = 1239
$ws.Range("K42").Value = 97
$ws.Range("K48").Value = 26
$ws.Range("K49").Value = 21
$ws.Range("K50").Value = 12
$ws.Range("K52").Value = 76
$ws.Range("I63").Value = 192
$ws.Range("J63").Value = 86
$ws.Range("K63").Value = 17
$ws.Range("K64").Value = 16
$ws.Range("K65").Value = 84
$ws.Range("K67").Value = 130
$ws.Range("K73").Value = 32
$ws.Range("J75").Value = 86
$ws.Range("K75").Value = 9
$ws.Range("K77").Value = 17
$ws.Range("K79").Value = 81
$ws.Range("K84").Value = 25
$ws.Range("K89").Value = 48
$ws.Range("K91").Value = 31
$ws.Range("K94").Value = 40
$ws.Range("K95").Value = 50
$ws.Range("K96").Value = 46
$ws.Range("K97").Value = 23
$ws.Range("I101").Value = 26237
$ws.Range("J101").Value = 29252
$ws.Range("K101").Value = 2990

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 37
$ws.Range("K3").Value = 22
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 21
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 24
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 5

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 8
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 34
$ws.Range("J4").Value = 57
$ws.Range("K4").Value = 4
$ws.Range("J7").Value = 1239
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 13
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 21

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 38
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J4").Value = 6
$ws.Range("K5").Value = 1
$ws.Range("J7").Value = 86
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 25
$ws.Range("K3").Value = 29
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 9
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 23
